$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns so values like "1.003" or "2.100"
# are kept as literal text (matching the source data) instead of being
# auto-coerced to numbers by Excel's smart entry.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.416.15"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "1.869.33"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "315.34"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "0.4664"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("D8").Value = "0.3731"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").Value = "0.07402"
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("D10").Value = "0.8915"
$ws.Range("E10").Value = "  +3.97%  "
$ws.Range("D11").Value = "0.07962"
$ws.Range("E11").Value = "  +5.95%  "
$ws.Range("D12").Value = "20.16"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").Value = "1.831.86"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "5.446"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "6.612"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "92.85"
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "0.000008960"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "14.95"
$ws.Range("E20").Value = "  +3.81%  "
$ws.Range("D21").Value = "27.425.76"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").Value = "5.173"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "2.036.47"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").Value = "152.58"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "1.867"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").Value = "18.61"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("D28").Value = "2.100"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").Value = "5.172"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").Value = "117.47"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").Value = "0.08924"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "0.7581"
$ws.Range("E32").Value = "  +5.85%  "
$ws.Range("D33").Value = "2.969"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").Value = "1.164"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").Value = "4.513"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("D36").Value = "2.589"
$ws.Range("E36").Value = "  +6.21%  "
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "0.05304"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "0.01959"
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("D40").Value = "2.994"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").Value = "7.166"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "0.5226"
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").Value = "0.1648"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").Value = "8.370"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D45").Value = "0.4919"
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("D46").Value = "10.41"
$ws.Range("E46").Value = "  +3.75%  "
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "103.42"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("D49").Value = "1.653"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("D50").Value = "0.06278"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "65.92"
$ws.Range("E51").Value = "  +3.12%  "
